# Updates cryptos list figures (price + 1h volume change) per the
# "Updated cryptos list on Fri Oct 13 09:09:04 UTC 2023 with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the value to be stored as text even when it looks like a number
    # (e.g. "206.93"), matching the original inlineStr cell type, by using a
    # leading apostrophe (Excel's text-quote prefix) and then restoring the
    # cell to the workbook default style so no stray formatting is introduced.
    $Range.Value = "'$Text"
    $Range.Style = 'Normal'
}

# Row 2
$ws.Range('D2').Value = '26.950.33'
$ws.Range('E2').Value = '  +0.43%  '

# Row 3
$ws.Range('D3').Value = '1.556.21'
$ws.Range('E3').Value = '  -0.24%  '

# Row 4
$ws.Range('E4').Value = '  +0.50%  '

# Row 5
Set-TextValue $ws.Range('D5') '206.93'
$ws.Range('E5').Value = '  +0.80%  '

# Row 6
$ws.Range('E6').Value = '  +1.15%  '

# Row 7
$ws.Range('E7').Value = '  +0.45%  '

# Row 8
$ws.Range('E8').Value = '  +0.33%  '

# Row 9
$ws.Range('E9').Value = '  +0.03%  '

# Row 10
$ws.Range('E10').Value = '  -0.24%  '

# Row 11
$ws.Range('E11').Value = '  -0.31%  '

# Row 12
$ws.Range('D12').Value = '1.777.97'
$ws.Range('E12').Value = '  -0.19%  '

# Row 13
$ws.Range('D13').Value = '1.559.07'
$ws.Range('E13').Value = '  +0.01%  '

# Row 15
$ws.Range('E15').Value = '  +0.29%  '

# Row 16
$ws.Range('D16').Value = '26.942.32'
$ws.Range('E16').Value = '  +0.40%  '

# Row 17
Set-TextValue $ws.Range('D17') '61.77'
$ws.Range('E17').Value = '  +0.75%  '

# Row 18
Set-TextValue $ws.Range('D18') '214.47'
$ws.Range('E18').Value = '  -0.04%  '

# Row 19
$ws.Range('E19').Value = '  +0.79%  '

# Row 21
$ws.Range('E21').Value = '  +0.48%  '

# Row 22
Set-TextValue $ws.Range('D22') '4.04'
$ws.Range('E22').Value = '  -1.86%  '

# Row 23
$ws.Range('E23').Value = '  +0.69%  '

# Row 24
$ws.Range('E24').Value = '  -2.60%  '

# Row 25
Set-TextValue $ws.Range('D25') '153.53'
$ws.Range('E25').Value = '  +0.03%  '

# Row 26
Set-TextValue $ws.Range('D26') '6.67'
$ws.Range('E26').Value = '  +0.87%  '

# Row 27
$ws.Range('E27').Value = '  -0.58%  '

# Row 29
$ws.Range('E29').Value = '  +0.76%  '

# Row 30
$ws.Range('E30').Value = '  -1.31%  '

# Row 31
$ws.Range('E31').Value = '  -0.60%  '

# Row 32
$ws.Range('E32').Value = '  +1.49%  '

# Row 33
$ws.Range('D33').Value = '1.369.76'
$ws.Range('E33').Value = '  -0.44%  '

# Row 34
$ws.Range('E34').Value = '  +1.47%  '

# Row 35
$ws.Range('E35').Value = '  +2.59%  '

# Row 36
Set-TextValue $ws.Range('D36') '0.971'
$ws.Range('E36').Value = '  +5.62%  '

# Row 37
$ws.Range('E37').Value = '  +0.52%  '

# Row 38
$ws.Range('E38').Value = '  +0.66%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.519'
$ws.Range('E39').Value = '  -1.06%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.807'
$ws.Range('E40').Value = '  -0.15%  '

# Row 41
$ws.Range('E41').Value = '  +0.43%  '

# Row 42
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D42') '0.981'
$ws.Range('E42').Value = '  -1.00%  '

# Row 43
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D43') '5.50'
$ws.Range('E43').Value = '  -1.08%  '

# Row 44
$ws.Range('E44').Value = '  +2.99%  '

# Row 45
Set-TextValue $ws.Range('D45') '63.71'
$ws.Range('E45').Value = '  +0.44%  '

# Row 46
Set-TextValue $ws.Range('D46') '1.73'
$ws.Range('E46').Value = '  -2.48%  '

# Row 47
$ws.Range('D47').Value = '1.690.69'
$ws.Range('E47').Value = '  -0.37%  '

# Row 48
Set-TextValue $ws.Range('D48') '86.06'
$ws.Range('E48').Value = '  -0.44%  '

# Row 49
$ws.Range('E49').Value = '  -1.03%  '

# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0975'
$ws.Range('E50').Value = '  -0.20%  '

# Row 51
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D51') '0.0956'
$ws.Range('E51').Value = '  +0.59%  '
